$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Introduction ")
$ws2 = $wb.Worksheets.Item("Cross review points ")

# ---------------------------------------------------------------
# Sheet "Introduction " updates
# ---------------------------------------------------------------

# Ref Version: 1 -> 1.4
$ws1.Range("D7").Value2 = 1.4

# Last update: text "24/01/2020" -> real date 2020-09-02 (serial 44076)
$ws1.Range("D9").Value2 = 44076
$ws1.Range("D9").NumberFormat = "m/d/yyyy"

# History table: new row for version 0.2
$ws1.Range("B14").Value2 = 0.2
$ws1.Range("C14").Value2 = "T.Sharaby"
$ws1.Range("E14").Value2 = 44076
$ws1.Range("E14").NumberFormat = "m/d/yyyy"
$ws1.Range("G14").Value2 = "Update the status and add one point "

# ---------------------------------------------------------------
# Sheet "Cross review points " updates
# ---------------------------------------------------------------

# Status column: Open -> Resolved for the first three open points
$ws2.Range("H2").Value2 = "Resolved"
$ws2.Range("H3").Value2 = "Resolved"
$ws2.Range("H4").Value2 = "Resolved"

# New open point / review row
$ws2.Range("A5").Value2 = 44076
$ws2.Range("A5").NumberFormat = "m/d/yyyy"
$ws2.Range("B5").Value2 = "T.Sharaby"
$ws2.Range("C5").Value2 = "There is a diffrenece between last modification date between the history table and status table "
$ws2.Range("D5").Value2 = "HSI_"
$ws2.Range("E5").Value2 = "ALL"
$ws2.Range("F5").Value2 = "Update the status and add one point "
$ws2.Range("H5").Value2 = "Open"
